$wb = $excel.ActiveWorkbook

# This script applies the numeric updates captured in the commit diff for
# "Siren_Profits" (the FFXIV leve-profit tracker spread across the ALC/ARM/
# BSM/CRP/CUL/GSM/LTW/WVR worksheets). Each worksheet row describes one Leve
# craft; columns H-N are derived market-price/profit figures that were
# refreshed by the scheduled data-collection run. No formulas are involved --
# every cell below is a literal numeric value in the source workbook.

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2457.7646
$ws.Range("I18").Value = 2392.9333
$ws.Range("J18").Value = 2944
$ws.Range("K18").Value = 2392.9333
$ws.Range("L18").Value = 2944
$ws.Range("M18").Value = -2108.9333
$ws.Range("N18").Value = -3512
$ws.Range("H74").Value = 6662.75
$ws.Range("I74").Value = 5160.4
$ws.Range("J74").Value = 9166.666999999999
$ws.Range("K74").Value = 5160.4
$ws.Range("L74").Value = 9166.666999999999
$ws.Range("M74").Value = -4224.4
$ws.Range("N74").Value = -11038.667
$ws.Range("H77").Value = 6662.75
$ws.Range("I77").Value = 5160.4
$ws.Range("J77").Value = 9166.666999999999
$ws.Range("K77").Value = 25802
$ws.Range("L77").Value = 45833.335
$ws.Range("M77").Value = -21122
$ws.Range("N77").Value = -55193.335
$ws.Range("H100").Value = 17074448
$ws.Range("I100").Value = 16672306
$ws.Range("J100").Value = 18280878
$ws.Range("K100").Value = 16672306
$ws.Range("L100").Value = 18280878
$ws.Range("M100").Value = -16671765
$ws.Range("N100").Value = -18281960
$ws.Range("H123").Value = 155555
$ws.Range("J123").Value = 155555
$ws.Range("L123").Value = 155555
$ws.Range("N123").Value = -165355
$ws.Range("H138").Value = 9160.450999999999
$ws.Range("I138").Value = 9897.700000000001
$ws.Range("J138").Value = 8809.380999999999
$ws.Range("K138").Value = 29693.1
$ws.Range("L138").Value = 26428.143
$ws.Range("M138").Value = -24553.1
$ws.Range("N138").Value = -36708.143

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 240941
$ws.Range("I45").Value = 359278.16
$ws.Range("K45").Value = 359278.16
$ws.Range("M45").Value = -358901.16

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3689.923
$ws.Range("I20").Value = 1829
$ws.Range("K20").Value = 1829
$ws.Range("M20").Value = -1582
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -32122
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -100608
$ws.Range("H94").Value = 9228.135
$ws.Range("I94").Value = 12166.64
$ws.Range("J94").Value = 3106.25
$ws.Range("K94").Value = 12166.64
$ws.Range("L94").Value = 3106.25
$ws.Range("M94").Value = -11715.64
$ws.Range("N94").Value = -4008.25
$ws.Range("H107").Value = 3630.2307
$ws.Range("I107").Value = 3572.0908
$ws.Range("J107").Value = 3950
$ws.Range("K107").Value = 3572.0908
$ws.Range("L107").Value = 3950
$ws.Range("M107").Value = -1652.0908
$ws.Range("N107").Value = -7790
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H130").Value = 99333
$ws.Range("J130").Value = 99333
$ws.Range("L130").Value = 99333
$ws.Range("N130").Value = -109373
$ws.Range("H134").Value = 10394.556
$ws.Range("I134").Value = 10928.25
$ws.Range("K134").Value = 32784.75
$ws.Range("M134").Value = -30249.75

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8078.12
$ws.Range("J58").Value = 3989.4167
$ws.Range("L58").Value = 3989.4167
$ws.Range("N58").Value = -4395.4167
$ws.Range("H86").Value = 17206.691
$ws.Range("I86").Value = 13983.333
$ws.Range("J86").Value = 19969.572
$ws.Range("K86").Value = 13983.333
$ws.Range("L86").Value = 19969.572
$ws.Range("M86").Value = -12860.333
$ws.Range("N86").Value = -22215.572
$ws.Range("H89").Value = 17206.691
$ws.Range("I89").Value = 13983.333
$ws.Range("J89").Value = 19969.572
$ws.Range("K89").Value = 69916.66500000001
$ws.Range("L89").Value = 99847.86
$ws.Range("M89").Value = -64300.66500000001
$ws.Range("N89").Value = -111079.86
$ws.Range("H132").Value = 33386206
$ws.Range("I132").Value = 41682584
$ws.Range("K132").Value = 125047752
$ws.Range("M132").Value = -125045222
$ws.Range("H136").Value = 8078.12
$ws.Range("J136").Value = 3989.4167
$ws.Range("L136").Value = 11968.2501
$ws.Range("N136").Value = -17068.2501

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 264004.84
$ws.Range("I5").Value = 534.1667
$ws.Range("J5").Value = 715668.9
$ws.Range("K5").Value = 1602.5001
$ws.Range("L5").Value = 2147006.7
$ws.Range("M5").Value = -1490.5001
$ws.Range("N5").Value = -2147230.7
$ws.Range("H114").Value = 8576
$ws.Range("I114").Value = 837.6667
$ws.Range("J114").Value = 15208.857
$ws.Range("K114").Value = 2513.0001
$ws.Range("L114").Value = 45626.571
$ws.Range("M114").Value = 740.9998999999998
$ws.Range("N114").Value = -52134.571
$ws.Range("H122").Value = 2495.9666
$ws.Range("I122").Value = 1434.3334
$ws.Range("J122").Value = 2613.926
$ws.Range("K122").Value = 12909.0006
$ws.Range("L122").Value = 23525.334
$ws.Range("M122").Value = -10459.0006
$ws.Range("N122").Value = -28425.334
$ws.Range("H135").Value = 264004.84
$ws.Range("I135").Value = 534.1667
$ws.Range("J135").Value = 715668.9
$ws.Range("K135").Value = 4807.5003
$ws.Range("L135").Value = 6441020.100000001
$ws.Range("M135").Value = -2272.5003
$ws.Range("N135").Value = -6446090.100000001

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6033.4585
$ws.Range("I70").Value = 6263.1055
$ws.Range("J70").Value = 5160.8
$ws.Range("K70").Value = 6263.1055
$ws.Range("L70").Value = 5160.8
$ws.Range("M70").Value = -5993.1055
$ws.Range("N70").Value = -5700.8
$ws.Range("H73").Value = 6033.4585
$ws.Range("I73").Value = 6263.1055
$ws.Range("J73").Value = 5160.8
$ws.Range("K73").Value = 6263.1055
$ws.Range("L73").Value = 5160.8
$ws.Range("M73").Value = -5327.1055
$ws.Range("N73").Value = -7032.8
$ws.Range("H80").Value = 5284.4585
$ws.Range("J80").Value = 3610.4614
$ws.Range("L80").Value = 3610.4614
$ws.Range("N80").Value = -5606.4614
$ws.Range("H83").Value = 5284.4585
$ws.Range("J83").Value = 3610.4614
$ws.Range("L83").Value = 18052.307
$ws.Range("N83").Value = -28036.307

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 25399.2
$ws.Range("I7").Value = 45687.6
$ws.Range("J7").Value = 5110.8
$ws.Range("K7").Value = 45687.6
$ws.Range("L7").Value = 5110.8
$ws.Range("M7").Value = -45575.6
$ws.Range("N7").Value = -5334.8
$ws.Range("H22").Value = 4303.803
$ws.Range("I22").Value = 6149.517
$ws.Range("J22").Value = 2631.125
$ws.Range("K22").Value = 6149.517
$ws.Range("L22").Value = 2631.125
$ws.Range("M22").Value = -5854.517
$ws.Range("N22").Value = -3221.125
$ws.Range("H27").Value = 4303.803
$ws.Range("I27").Value = 6149.517
$ws.Range("J27").Value = 2631.125
$ws.Range("K27").Value = 6149.517
$ws.Range("L27").Value = 2631.125
$ws.Range("M27").Value = -6042.517
$ws.Range("N27").Value = -2845.125
$ws.Range("H55").Value = 588.4054
$ws.Range("I55").Value = 611.4286
$ws.Range("K55").Value = 611.4286
$ws.Range("M55").Value = -438.4286
$ws.Range("H61").Value = 2341.8096
$ws.Range("I61").Value = 1880.5
$ws.Range("J61").Value = 2761.182
$ws.Range("K61").Value = 1880.5
$ws.Range("L61").Value = 2761.182
$ws.Range("M61").Value = -1678.5
$ws.Range("N61").Value = -3165.182
$ws.Range("H68").Value = 4408.7856
$ws.Range("I68").Value = 2749.3333
$ws.Range("K68").Value = 2749.3333
$ws.Range("M68").Value = -2000.3333
$ws.Range("H71").Value = 4408.7856
$ws.Range("I71").Value = 2749.3333
$ws.Range("K71").Value = 13746.6665
$ws.Range("M71").Value = -10002.6665
$ws.Range("H82").Value = 2793.3333
$ws.Range("J82").Value = 2293.6667
$ws.Range("L82").Value = 2293.6667
$ws.Range("N82").Value = -3015.6667
$ws.Range("H85").Value = 2793.3333
$ws.Range("J85").Value = 2293.6667
$ws.Range("L85").Value = 2293.6667
$ws.Range("N85").Value = -4789.6667
$ws.Range("H113").Value = 2341.8096
$ws.Range("I113").Value = 1880.5
$ws.Range("J113").Value = 2761.182
$ws.Range("K113").Value = 1880.5
$ws.Range("L113").Value = 2761.182
$ws.Range("M113").Value = 289.5
$ws.Range("N113").Value = -7101.182
$ws.Range("H126").Value = 25399.2
$ws.Range("I126").Value = 45687.6
$ws.Range("J126").Value = 5110.8
$ws.Range("K126").Value = 137062.8
$ws.Range("L126").Value = 15332.4
$ws.Range("M126").Value = -134592.8
$ws.Range("N126").Value = -20272.4

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11951.6
$ws.Range("J41").Value = 11279.556
$ws.Range("L41").Value = 11279.556
$ws.Range("N41").Value = -12059.556
$ws.Range("H62").Value = 137138.44
$ws.Range("I62").Value = 268999.44
$ws.Range("K62").Value = 268999.44
$ws.Range("M62").Value = -268375.44
$ws.Range("H65").Value = 137138.44
$ws.Range("I65").Value = 268999.44
$ws.Range("K65").Value = 1344997.2
$ws.Range("M65").Value = -1341877.2
$ws.Range("H96").Value = 39395384
$ws.Range("I96").Value = 14287199
$ws.Range("K96").Value = 14287199
$ws.Range("M96").Value = -14285826
$ws.Range("H140").Value = 81044.60000000001
$ws.Range("J140").Value = 81044.60000000001
$ws.Range("L140").Value = 81044.60000000001
$ws.Range("N140").Value = -91404.60000000001
